$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 currently holds phone "09876543" (text) with 0 points.
# A new entry for the same phone number earns 120 points, so duplicate
# row 31 into a new row 32 (preserving its text formatting / leading
# zero) and set its points to 120. The original row 31's phone value
# is stored as a plain number going forward.

$ws.Rows("32").Insert()
$ws.Range("A31:C31").Copy($ws.Range("A32:C32"))

$ws.Cells.Item(32, 3).Value = 120
$ws.Cells.Item(31, 1).Value = 9876543
